$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.321162700653076
$ws.Range("B1").Value = 2.163686037063599
$ws.Range("C1").Value = 2.772473812103271
$ws.Range("D1").Value = 2.68518328666687
$ws.Range("E1").Value = 0.9968612790107727
